# Switched out the wrong bearing link for the shoulder rot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://www.amazon.com/Timken-SET37-Tapered-Roller-Bearing/dp/B000BZ6YEK?crid=2UVMW157QN9VD&dib=eyJ2IjoiMSJ9.2iRD0--LwAOtK6sCHU59ccwGJuRc5xFBYyjJQh18M4UXnO7H_ZEp_RF1DBAneDyUPEW41QxDKmbDABDOTskGcetuw5M7rczrcJk2ijXBzWPzPx4e81ycuzjMNTUf1LCgl8F85OwKQvS0FSy3HqugdxaoWSM9lra-DHRZzyeZ4RZV037QmoTEjlaHqrlmwPh9t85RVao1pSw0P0kDZEYdu48P_FysqfBYEduzfz9srAsN7hSKuiP0O7mAwPilwJhbmhWH6cM_JCfLtP8VQdRd-V7lYpPKFxGmKKUMIfvId7g.83HkNvFnLf5ILJIPVbD65pyGfEMhvFL6FLQp0veUnws&dib_tag=se&keywords=tapered+roller+bearing&qid=1751500342&s=industrial&sprefix=tapered+roller+beairn%2Cindustrial%2C164&sr=1-6"
$newName = "Shoulder Rot Bearing - SET37 * 2"

# Update the link cell text first so the new URL string is allocated
# ahead of the new part-name string in the shared-string table (matches
# the order Excel itself would have produced).
$ws.Range("C9").Value = $newUrl
$ws.Range("A9").Value = $newName
$ws.Range("D9").Value = 35.36

# Re-point the existing hyperlink on C9 to the new product page.
$linkCell = $null
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq "`$C`$9") {
        $linkCell = $h
    }
}

if ($linkCell -ne $null) {
    $linkCell.Address = $newUrl
    $linkCell.TextToDisplay = $newUrl
}

# Move the active selection, matching the saved workbook view.
$ws.Range("E11").Select()
